# Generate Report for Handoff
#
# The "aafca9a8-6ae7-489c-8d3b-60c15ee3eb67.md" file has finished its
# localization round for zh-cn and de-de: it is no longer "Handed back: in
# sync with en-US" but is now "Ready for handoff" again, the handoff/handback
# timestamps move forward, and the zh-cn localization was produced via
# machine translation ("mt" instead of "ht"). The handback commit found for
# this file isn't the latest revision on GitHub, so an Error Detail message
# is recorded for it in both language tabs.

$wb = $excel.ActiveWorkbook

$errorDetail = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9ea5e930fa5eb4e7fd99b334c2fa8f12b2bc00eb/e2e/aafca9a8-6ae7-489c-8d3b-60c15ee3eb67.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d0418897c41ce1a4b28b9d2e40f4df2ad7b30786/e2e/aafca9a8-6ae7-489c-8d3b-60c15ee3eb67.md."

# Column-width constants. Excel's ColumnWidth property is expressed in
# "characters" and gets snapped to the sheet's pixel grid, so the values
# below are chosen so the round-tripped width matches the target display
# width used by the handoff-report generator (~17.22 chars for the
# shrunken Status/Latest Handoff columns, 40 chars for the widened
# Error Detail column).
$narrowColWidth = 16 + 1/3
$wideColWidth = 39 + 1/6

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-10-20 01:18:31"

$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-10-20 01:18:31"

$wsOverview.Columns.Item(5).ColumnWidth = $narrowColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowColWidth

# --- zh-cn sheet -----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-10-20 01:18:20"

$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-10-20 01:18:20"
$wsZhCn.Range("P3").Value = $errorDetail

$wsZhCn.Columns.Item(3).ColumnWidth = $narrowColWidth
$wsZhCn.Columns.Item(16).ColumnWidth = $wideColWidth

# --- de-de sheet -----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-10-20 01:18:31"

$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-10-20 01:18:31"
$wsDeDe.Range("P3").Value = $errorDetail

$wsDeDe.Columns.Item(3).ColumnWidth = $narrowColWidth
$wsDeDe.Columns.Item(16).ColumnWidth = $wideColWidth
